$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Benefit" -> "Status"
$ws.Range("D1").Value = "Status"

# Column D now holds a win(1)/loss(0) flag instead of the computed benefit amount.
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 1
$ws.Range("D15").Value = 0

# Update the active selection to D3
$ws.Range("D3").Select()
